# Update Soliera.xlsx with new daily records through 23 agosto 2021
# Adds rows 344-357 (dates 44418-44431) to the existing data table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(344, 44418, 0, 16, 103.3391461603049),
    @(345, 44419, 0, 15, 96.8804495252858),
    @(346, 44420, 2, 12, 77.50435962022863),
    @(347, 44421, 1, 10, 64.58696635019054),
    @(348, 44422, 9, 17, 109.7978427953239),
    @(349, 44423, 3, 17, 109.7978427953239),
    @(350, 44424, 7, 22, 142.0913259704192),
    @(351, 44425, 1, 23, 148.5500226054382),
    @(352, 44426, 2, 25, 161.4674158754763),
    @(353, 44427, 0, 23, 148.5500226054382),
    @(354, 44428, 4, 26, 167.9261125104954),
    @(355, 44429, 5, 22, 142.0913259704192),
    @(356, 44430, 0, 19, 122.715236065362),
    @(357, 44431, 1, 13, 83.96305625524769)
)

# Reference cell A343 carries the date style used for column A
$styleSourceA = $ws.Range("A343")

foreach ($row in $data) {
    $r = $row[0]

    $cellA = $ws.Cells.Item($r, 1)
    $cellA.Value = $row[1]
    $styleSourceA.Copy()
    $cellA.PasteSpecial(-4122) | Out-Null

    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}

$excel.CutCopyMode = 0
